# Swap the order of "System" and the email address in the
# "Recorded By" column (column G) wherever the combined value
# "System, dnasr281@gmail.com" appears, turning it into
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Column G holds the "Recorded By" values.
$colIndex = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
